$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.735.89"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "3.125.90"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.47"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.62"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "3.124.43"
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.476"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +6.18%  "
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.107"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.413"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.71%  "
$ws.Range("E13").Value = "  +1.54%  "
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.91"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("D17").Value = "57.844.59"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "3.130.60"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("E19").Value = "  +2.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.72"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.07"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "367.74"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +6.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("E24").Value = "  -2.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.23"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.33%  "
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").Value = "0.0₃0863"
$ws.Range("E29").Value = "  -2.72%  "
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.07"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.43"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.15"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.17"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.78"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("E38").Value = "  +5.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.40"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.57%  "
$ws.Range("E40").Value = "  +4.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0673"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.48%  "
$ws.Range("D42").Value = "2.533.33"
$ws.Range("E42").Value = "  +6.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.09"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "37.69"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0269"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.976"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("E49").Value = "  +2.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.77"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("E51").Value = "  -1.75%  "
